$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUuid = "05adc59e-9ff8-4855-a157-5b5f0cf7ac22"

# --- Update existing rows 2-5: refresh the uuid value in column G ---
$ws.Range("G2").Value = $newUuid
$ws.Range("G3").Value = $newUuid
$ws.Range("G4").Value = $newUuid
$ws.Range("G5").Value = $newUuid

# --- Remove the now-unused numeric style from the downtime column (H2:H5) ---
$ws.Range("H2").Style = "Normal"
$ws.Range("H3").Style = "Normal"
$ws.Range("H4").Style = "Normal"
$ws.Range("H5").Style = "Normal"

# --- Append two new data rows (6 and 7) from the Highcharts visualization export ---

# Row 6 (mirrors row 2's "pri cl LA" reading)
$ws.Range("A6").Value = "Line:8 Stage:1"
$ws.Range("B6").Value = "'01/09/2024"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "pri cl LA"
$ws.Range("D6").Value = 45300.42056299769
$ws.Range("D6").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E6").Value = 45300.42067873842
$ws.Range("E6").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("F6").Value = 0.17
$ws.Range("G6").Value = $newUuid
$ws.Range("H6").Value = 10

# Row 7 (mirrors row 3's "pri pH HA" reading)
$ws.Range("A7").Value = "Line:8 Stage:1"
$ws.Range("B7").Value = "'01/09/2024"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "pri pH HA"
$ws.Range("D7").Value = 45300.65476473379
$ws.Range("D7").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E7").Value = 45300.65488047454
$ws.Range("E7").NumberFormat = $ws.Range("E3").NumberFormat
$ws.Range("F7").Value = 0.17
$ws.Range("G7").Value = $newUuid
$ws.Range("H7").Value = 10
